# Updated symbol list on Tue Dec 13 14:46:12 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
}

Set-TextValue $ws.Range("D2") "278.45"
Set-TextValue $ws.Range("D3") "22.80"
Set-TextValue $ws.Range("D4") "6.364"
Set-TextValue $ws.Range("D5") "0.06273"
Set-TextValue $ws.Range("D6") "3.666"
Set-TextValue $ws.Range("D7") "6.624"
Set-TextValue $ws.Range("D8") "1.393"
Set-TextValue $ws.Range("D9") "0.8302"
Set-TextValue $ws.Range("D10") "0.01383"
Set-TextValue $ws.Range("D12") "0.08405"
Set-TextValue $ws.Range("D13") "0.03502"
Set-TextValue $ws.Range("D14") "0.03227"
Set-TextValue $ws.Range("D15") "4.094"
Set-TextValue $ws.Range("D16") "0.09300"
Set-TextValue $ws.Range("D17") "0.001681"
Set-TextValue $ws.Range("D18") "0.04755"
Set-TextValue $ws.Range("D19") "0.006256"
Set-TextValue $ws.Range("D20") "0.005731"
Set-TextValue $ws.Range("D21") "0.001077"
Set-TextValue $ws.Range("D22") "0.0001497"
Set-TextValue $ws.Range("D23") "3.727"
Set-TextValue $ws.Range("D25") "0.3333"
Set-TextValue $ws.Range("D26") "0.1260"
Set-TextValue $ws.Range("D28") "0.0002701"
Set-TextValue $ws.Range("D40") "0.04743"
Set-TextValue $ws.Range("D41") "0.007096"
Set-TextValue $ws.Range("D43") "0.003442"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.01225"
Set-TextValue $ws.Range("D45") "0.00006090"
Set-TextValue $ws.Range("D46") "0.0009888"
Set-TextValue $ws.Range("D48") "0.7811"
Set-TextValue $ws.Range("D49") "0.002475"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"
